# Update odds values in row 2 of the active worksheet to reflect the
# latest FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 2.5
    "H2"  = 2.9
    "M2"  = 1.13
    "N2"  = 6
    "O2"  = 1.53
    "P2"  = 2.38
    "Q2"  = 2.7
    "R2"  = 1.44
    "S2"  = 1.57
    "U2"  = 2.2
    "V2"  = 1.62
    "X2"  = 11
    "Y2"  = 11
    "Z2"  = 26
    "AA2" = 26
    "AC2" = 6
    "AF2" = 81
    "AG2" = 7
    "AP2" = 34
    "AR2" = 101
    "AT2" = 2.2
    "AU2" = 9.5
    "AW2" = 4.75
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
